# Updated cryptos list on Mon May 15 08:16:27 UTC 2023 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) figures for each coin row
# on Sheet1, and fixes two rows where the coin order/name+link had been swapped
# (Frax/TrustWalletToken at rows 36-37, Decentraland/NEARProtocol at rows 47-48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.801.07"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.854.11"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.016"
$ws.Range("E4").Value = "  -2.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.03"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4328"
$ws.Range("E7").Value = "  -2.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3787"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07415"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8853"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.74"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "1.858.96"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.763"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.495"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07166"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.68"
$ws.Range("E16").Value = "  +5.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.018"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009047"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.013"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "27.768.98"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.280"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").Value = "2.092.44"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.038"
$ws.Range("E25").Value = "  +3.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.37"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.69"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.067"
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.04"
$ws.Range("E30").Value = "  +3.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08974"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.241"
$ws.Range("E32").Value = "  +2.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7809"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.931"
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.147"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.015"
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05333"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01970"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.885"
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5204"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.034"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1686"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.846"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "111.01"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.80"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.714"
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4749"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06529"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.016"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.896"
$ws.Range("E51").Value = "  +0.62%  "

Write-Host "Applied cryptos update"
